# Final commit by Shilpi
# Applies the data edits + selection/active-tab changes captured in the diff.

$wb = $excel.ActiveWorkbook

$wsPostUsers = $wb.Worksheets.Item("PostUsers")
$wsExistOneField = $wb.Worksheets.Item("PostUsersExistOneField")

# --- Sheet "PostUsers": rename the first_name values in column A (rows 2-6) ---
$wsPostUsers.Range("A2").Value = "Kirtii"
$wsPostUsers.Range("A3").Value = "Shaneew"
$wsPostUsers.Range("A4").Value = "Shaunre"
$wsPostUsers.Range("A5").Value = "Duncantt"
$wsPostUsers.Range("A6").Value = "Priyate"

# New (until-now unused) column L touched on every data row - stays blank,
# matching the trailing empty <c r="Lx"/> cells that show up after a resave.
$wsPostUsers.Range("L2").ClearContents()
$wsPostUsers.Range("L2").Style = "Normal"
$wsPostUsers.Range("L3").ClearContents()
$wsPostUsers.Range("L3").Style = "Normal"
$wsPostUsers.Range("L4").ClearContents()
$wsPostUsers.Range("L4").Style = "Normal"
$wsPostUsers.Range("L5").ClearContents()
$wsPostUsers.Range("L5").Style = "Normal"
$wsPostUsers.Range("L6").ClearContents()
$wsPostUsers.Range("L6").Style = "Normal"

# --- Sheet "PostUsersExistOneField": data corrections ---
$wsExistOneField.Range("A2").Value = "Kiran"
$wsExistOneField.Range("D3").Value = 1234097089
$wsExistOneField.Range("E4").Value = "abct@xyz.com"

# Column L values were removed on this sheet - clear content but keep the cells present.
$wsExistOneField.Range("L2").ClearContents()
$wsExistOneField.Range("L2").Style = "Normal"
$wsExistOneField.Range("L3").ClearContents()
$wsExistOneField.Range("L3").Style = "Normal"
$wsExistOneField.Range("L4").ClearContents()
$wsExistOneField.Range("L4").Style = "Normal"

# --- Selection / active sheet changes ---
# End the session with "PostUsers" showing A6 selected (no longer the active tab)...
$wsPostUsers.Range("A6").Select()
# ...and "PostUsersExistOneField" as the active tab with E4 selected.
$wsExistOneField.Activate()
$wsExistOneField.Range("E4").Select()
